$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.637.63'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '3.118.65'
$ws.Range("E3").Value = '  +0.43%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '531.42'
$ws.Range("E5").Value = '  +1.41%  '
$ws.Range("D6").Value = '138.06'
$ws.Range("E6").Value = '  +1.32%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '3.117.86'
$ws.Range("E8").Value = '  +0.46%  '
$ws.Range("D9").Value = '''0.470'
$ws.Range("E9").Value = '  +5.35%  '
$ws.Range("E10").Value = '  +0.30%  '
$ws.Range("E11").Value = '  +0.55%  '
$ws.Range("D12").Value = '0.411'
$ws.Range("E12").Value = '  +4.52%  '
$ws.Range("E13").Value = '  +1.40%  '
$ws.Range("D14").Value = '3.653.28'
$ws.Range("E14").Value = '  +0.33%  '
$ws.Range("D15").Value = '25.61'
$ws.Range("E15").Value = '  +1.75%  '
$ws.Range("E16").Value = '  +1.11%  '
$ws.Range("D17").Value = '57.735.75'
$ws.Range("E17").Value = '  +0.50%  '
$ws.Range("D18").Value = '3.115.70'
$ws.Range("E18").Value = '  +0.33%  '
$ws.Range("D20").Value = '12.66'
$ws.Range("E20").Value = '  +1.97%  '
$ws.Range("D21").Value = '8.06'
$ws.Range("E21").Value = '  +2.72%  '
$ws.Range("D22").Value = '360.93'
$ws.Range("E22").Value = '  +4.26%  '
$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").Value = '69.07'
$ws.Range("E24").Value = '  +2.19%  '
$ws.Range("E25").Value = '  +0.72%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("D28").Value = '0.0₃0866'
$ws.Range("E28").Value = '  -3.07%  '
$ws.Range("D29").Value = '7.28'
$ws.Range("E29").Value = '  -1.18%  '
$ws.Range("D30").Value = '6.09'
$ws.Range("E30").Value = '  +0.80%  '
$ws.Range("D31").Value = '1.87'
$ws.Range("E31").Value = '  +0.17%  '
$ws.Range("D32").Value = '''21.30'
$ws.Range("E32").Value = '  +2.18%  '
$ws.Range("D33").Value = '5.14'
$ws.Range("E33").Value = '  +3.92%  '
$ws.Range("E34").Value = '  -0.34%  '
$ws.Range("D35").Value = '159.13'
$ws.Range("E35").Value = '  +0.68%  '
$ws.Range("D36").Value = '6.05'
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("E37").Value = '  +5.06%  '
$ws.Range("D38").Value = '25.46'
$ws.Range("E38").Value = '  -1.09%  '
$ws.Range("D39").Value = '1.67'
$ws.Range("E39").Value = '  +4.15%  '
$ws.Range("D40").Value = '0.0669'
$ws.Range("E40").Value = '  +1.30%  '
$ws.Range("D41").Value = '2.495.64'
$ws.Range("E41").Value = '  +5.96%  '
$ws.Range("D42").Value = '4.01'
$ws.Range("E42").Value = '  -3.51%  '
$ws.Range("E43").Value = '  -0.34%  '
$ws.Range("D44").Value = '37.79'
$ws.Range("E44").Value = '  +3.20%  '
$ws.Range("E45").Value = '  +1.27%  '
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  -0.05%  '
$ws.Range("E47").Value = '  +1.98%  '
$ws.Range("D48").Value = '6.09'
$ws.Range("E48").Value = '  +1.84%  '
$ws.Range("D49").Value = '19.66'
$ws.Range("E49").Value = '  -0.54%  '
$ws.Range("D50").Value = '0.741'
$ws.Range("E50").Value = '  -1.96%  '
$ws.Range("E51").Value = '  +2.47%  '
